$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Uncomment RAD Extension Payments test row: switch Execute flag from
# "DONOTRUN" back to "Y" so the test case actually runs.
$ws.Range("C4").Value = "Y"

# Leave the edited cell selected, matching the author's saved selection.
[void]$ws.Range("C4").Select()
